$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the single data row (row 2) with the new spot-price values
$ws.Range("A2").Value = 46001

$ws.Range("B2").Value = 97.13
$ws.Range("C2").Value = 87.48
$ws.Range("D2").Value = 85.45999999999999
$ws.Range("E2").Value = 81.8
$ws.Range("F2").Value = 84.09999999999999
$ws.Range("G2").Value = 87.43000000000001
$ws.Range("H2").Value = 96.48999999999999
$ws.Range("I2").Value = 108.25
$ws.Range("J2").Value = 114.39
$ws.Range("K2").Value = 108.52
$ws.Range("L2").Value = 104.27
$ws.Range("M2").Value = 100.65
$ws.Range("N2").Value = 98.43000000000001
$ws.Range("O2").Value = 97.42
$ws.Range("P2").Value = 98.45
$ws.Range("Q2").Value = 101.08
$ws.Range("R2").Value = 108.4
$ws.Range("S2").Value = 115.18
$ws.Range("T2").Value = 118.28
$ws.Range("U2").Value = 119.93
$ws.Range("V2").Value = 124.67
$ws.Range("W2").Value = 119.92
$ws.Range("X2").Value = 112.21
$ws.Range("Y2").Value = 102.22
$ws.Range("Z2").Value = 103.01

$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 115.45
$ws.Range("AD2").Value = 122.3
$ws.Range("AF2").Value = 119.1
$ws.Range("AG2").Value = "0h-23h"
